$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.829496891513515
$ws.Range("D2").Value = 7.807266741131357
$ws.Range("E2").Value = 12.84878488013075
$ws.Range("F2").Value = 38.55121716906317
$ws.Range("G2").Value = 44.5519243583317
$ws.Range("H2").Value = 18.1106718550022
$ws.Range("J2").Value = 10.10306241195581
$ws.Range("K2").Value = 13.55956539087151
$ws.Range("M2").Value = 16.91879994382163
$ws.Range("N2").Value = 20.54998897460586
$ws.Range("B3").Value = 7.759531012027722
$ws.Range("D3").Value = 7.799845262515726
$ws.Range("E3").Value = 12.85954614347802
$ws.Range("F3").Value = 38.54655489067291
$ws.Range("G3").Value = 44.47163725919039
$ws.Range("H3").Value = 18.14631825128028
$ws.Range("J3").Value = 10.12477245294603
$ws.Range("K3").Value = 13.29821569749464
$ws.Range("M3").Value = 16.83116551171566
$ws.Range("N3").Value = 20.61813602253495
$ws.Range("B4").Value = 7.718118484130859
$ws.Range("D4").Value = 7.796234140675544
$ws.Range("E4").Value = 12.86806121040808
$ws.Range("F4").Value = 38.55308244084056
$ws.Range("G4").Value = 44.43547865316583
$ws.Range("H4").Value = 18.17158333852938
$ws.Range("J4").Value = 10.13920507039803
$ws.Range("K4").Value = 13.138181254169
$ws.Range("M4").Value = 16.78037883314383
$ws.Range("N4").Value = 20.66189904060725
$ws.Range("B5").Value = 7.701649965198593
$ws.Range("D5").Value = 7.795001770148109
$ws.Range("E5").Value = 12.87201122921828
$ws.Range("F5").Value = 38.55810279248468
$ws.Range("G5").Value = 44.42405339397921
$ws.Range("H5").Value = 18.18272711173535
$ws.Range("J5").Value = 10.14536401913666
$ws.Range("K5").Value = 13.07316791292075
$ws.Range("M5").Value = 16.76045877407498
$ws.Range("N5").Value = 20.6802172196804
$ws.Range("B6").Value = 7.698940512139804
$ws.Range("D6").Value = 7.794811621091373
$ws.Range("E6").Value = 12.87269613026118
$ws.Range("F6").Value = 38.55907892025787
$ws.Range("G6").Value = 44.42235624082122
$ws.Range("H6").Value = 18.18462870331501
$ws.Range("J6").Value = 10.14640347915501
$ws.Range("K6").Value = 13.06238742585115
$ws.Range("M6").Value = 16.75719838025646
$ws.Range("N6").Value = 20.68328823906543
$ws.Range("B7").Value = 7.717894710449968
$ws.Range("D7").Value = 7.796216550252526
$ws.Range("E7").Value = 12.86811253757275
$ws.Range("F7").Value = 38.55314059257982
$ws.Range("G7").Value = 44.43531116218764
$ws.Range("H7").Value = 18.17173019562829
$ws.Range("J7").Value = 10.13928700800184
$ws.Range("K7").Value = 13.13730352046941
$ws.Range("M7").Value = 16.78010702196158
$ws.Range("N7").Value = 20.66214412274068
$ws.Range("B8").Value = 7.80506292094047
$ws.Range("D8").Value = 7.804512396510349
$ws.Range("E8").Value = 12.85209969335688
$ws.Range("F8").Value = 38.54766054470999
$ws.Range("G8").Value = 44.52151922724986
$ws.Range("H8").Value = 18.12226081439524
$ws.Range("J8").Value = 10.11031937226153
$ws.Range("K8").Value = 13.46941672107073
$ws.Range("M8").Value = 16.8879665821401
$ws.Range("N8").Value = 20.57308838161759
$ws.Range("B9").Value = 7.987427432585218
$ws.Range("D9").Value = 7.828219979258097
$ws.Range("E9").Value = 12.83581522777338
$ws.Range("F9").Value = 38.61140366070799
$ws.Range("G9").Value = 44.79439624444979
$ws.Range("H9").Value = 18.05211326206515
$ws.Range("J9").Value = 10.0622493992597
$ws.Range("K9").Value = 14.12016703637635
$ws.Range("M9").Value = 17.1226915378806
$ws.Range("N9").Value = 20.41361949311607
$ws.Range("B10").Value = 8.127205991997347
$ws.Range("D10").Value = 7.850082372516799
$ws.Range("E10").Value = 12.83303635628736
$ws.Range("F10").Value = 38.70353665969108
$ws.Range("G10").Value = 45.05738158673579
$ws.Range("H10").Value = 18.01702536370973
$ws.Range("J10").Value = 10.03224030709416
$ws.Range("K10").Value = 14.5928074019517
$ws.Range("M10").Value = 17.30821521214678
$ws.Range("N10").Value = 20.30560710286491
$ws.Range("B11").Value = 8.191796378730208
$ws.Range("D11").Value = 7.860971811980322
$ws.Range("E11").Value = 12.83375805395153
$ws.Range("F11").Value = 38.75523180053254
$ws.Range("G11").Value = 45.19035017317876
$ws.Range("H11").Value = 18.00464777056677
$ws.Range("J11").Value = 10.01973726755468
$ws.Range("K11").Value = 14.80561298096385
$ws.Range("M11").Value = 17.39520558508057
$ws.Range("N11").Value = 20.25843519578416
$ws.Range("B12").Value = 8.216377725111933
$ws.Range("D12").Value = 7.865229121018112
$ws.Range("E12").Value = 12.83431594045627
$ws.Range("F12").Value = 38.77620683586462
$ws.Range("G12").Value = 45.24259352818417
$ws.Range("H12").Value = 18.00047685124493
$ws.Range("J12").Value = 10.01516752254066
$ws.Range("K12").Value = 14.88580051400896
$ws.Range("M12").Value = 17.428497745881
$ws.Range("N12").Value = 20.24085321573499
$ws.Range("B13").Value = 8.211078601107543
$ws.Range("D13").Value = 7.864306320054329
$ws.Range("E13").Value = 12.83418314880669
$ws.Range("F13").Value = 38.7716273949646
$ws.Range("G13").Value = 45.23125834074686
$ws.Range("H13").Value = 18.00135216363326
$ws.Range("J13").Value = 10.01614436930063
$ws.Range("K13").Value = 14.86854960420601
$ws.Range("M13").Value = 17.42131244319784
$ws.Range("N13").Value = 20.24462733361757
$ws.Range("B14").Value = 8.193816374775167
$ws.Range("D14").Value = 7.861319395612369
$ws.Range("E14").Value = 12.83379825384883
$ws.Range("F14").Value = 38.75692943646813
$ws.Range("G14").Value = 45.19461053018743
$ws.Range("H14").Value = 18.00429427687298
$ws.Range("J14").Value = 10.01935800867768
$ws.Range("K14").Value = 14.81221846751297
$ws.Range("M14").Value = 17.39793765064423
$ws.Range("N14").Value = 20.25698309177702
$ws.Range("B15").Value = 8.18325804401575
$ws.Range("D15").Value = 7.859507169196693
$ws.Range("E15").Value = 12.83359952660431
$ws.Range("F15").Value = 38.74810847700772
$ws.Range("G15").Value = 45.17240811179201
$ws.Range("H15").Value = 18.00616365301207
$ws.Range("J15").Value = 10.02134792155123
$ws.Range("K15").Value = 14.77765995429272
$ws.Range("M15").Value = 17.38366492668554
$ws.Range("N15").Value = 20.2645879047686
$ws.Range("B16").Value = 8.123003183777156
$ws.Range("D16").Value = 7.849389547380882
$ws.Range("E16").Value = 12.83302907187462
$ws.Range("F16").Value = 38.70035451131518
$ws.Range("G16").Value = 45.04895792278876
$ws.Range("H16").Value = 18.01790646326035
$ws.Range("J16").Value = 10.03308049599832
$ws.Range("K16").Value = 14.57884875758503
$ws.Range("M16").Value = 17.30258050441063
$ws.Range("N16").Value = 20.30872929255812
$ws.Range("B17").Value = 8.086280342345036
$ws.Range("D17").Value = 7.843423101474074
$ws.Range("E17").Value = 12.83318708595332
$ws.Range("F17").Value = 38.67355983488178
$ws.Range("G17").Value = 44.9766228832015
$ws.Range("H17").Value = 18.02602890112709
$ws.Range("J17").Value = 10.04057197929768
$ws.Range("K17").Value = 14.45626190764618
$ws.Range("M17").Value = 17.2534867174126
$ws.Range("N17").Value = 20.33631053411582
$ws.Range("B18").Value = 8.065254201197847
$ws.Range("D18").Value = 7.840080354446459
$ws.Range("E18").Value = 12.83346495163462
$ws.Range("F18").Value = 38.65906966139854
$ws.Range("G18").Value = 44.93627466380784
$ws.Range("H18").Value = 18.03103802326279
$ws.Range("J18").Value = 10.04498897468371
$ws.Range("K18").Value = 14.38555059645225
$ws.Range("M18").Value = 17.22549487411841
$ws.Range("N18").Value = 20.35235942571053
$ws.Range("B19").Value = 8.058152248210337
$ws.Range("D19").Value = 7.838963903480558
$ws.Range("E19").Value = 12.83359117355086
$ws.Range("F19").Value = 38.65432199815189
$ws.Range("G19").Value = 44.92283009188656
$ws.Range("H19").Value = 18.0327919316254
$ws.Range("J19").Value = 10.04650306544526
$ws.Range("K19").Value = 14.36157670095502
$ws.Range("M19").Value = 17.21606017014466
$ws.Range("N19").Value = 20.35782510110818
$ws.Range("B20").Value = 8.090179789569319
$ws.Range("D20").Value = 7.844049043729255
$ws.Range("E20").Value = 12.83315091993443
$ws.Range("F20").Value = 38.67631686101382
$ws.Range("G20").Value = 44.98419316031775
$ws.Range("H20").Value = 18.0251293368293
$ws.Range("J20").Value = 10.03976331263489
$ws.Range("K20").Value = 14.46933305486868
$ws.Range("M20").Value = 17.25868757552544
$ws.Range("N20").Value = 20.33335533811843
$ws.Range("B21").Value = 8.198883563817375
$ws.Range("D21").Value = 7.862193115580402
$ws.Range("E21").Value = 12.83390359127734
$ws.Range("F21").Value = 38.76120867775757
$ws.Range("G21").Value = 45.20532378644398
$ws.Range("H21").Value = 18.00341609148782
$ws.Range("J21").Value = 10.01840961182038
$ws.Range("K21").Value = 14.82877569295754
$ws.Range("M21").Value = 17.40479405411343
$ws.Range("N21").Value = 20.25334629059063
$ws.Range("B22").Value = 8.270630824049006
$ws.Range("D22").Value = 7.874829685723606
$ws.Range("E22").Value = 12.83605381119366
$ws.Range("F22").Value = 38.82484138182688
$ws.Range("G22").Value = 45.36085320466478
$ws.Range("H22").Value = 17.99223433791227
$ws.Range("J22").Value = 10.00541470005415
$ws.Range("K22").Value = 15.06134206194289
$ws.Range("M22").Value = 17.50231704578313
$ws.Range("N22").Value = 20.20269300335737
$ws.Range("B23").Value = 8.232280911021705
$ws.Range("D23").Value = 7.868014792377785
$ws.Range("E23").Value = 12.83475481311991
$ws.Range("F23").Value = 38.79013652413808
$ws.Range("G23").Value = 45.27684664252718
$ws.Range("H23").Value = 17.99792668929492
$ws.Range("J23").Value = 10.01226247735902
$ws.Range("K23").Value = 14.93745758956079
$ws.Range("M23").Value = 17.45008860004906
$ws.Range("N23").Value = 20.2295782498588
$ws.Range("B24").Value = 8.088416578962773
$ws.Range("D24").Value = 7.843765782724656
$ws.Range("E24").Value = 12.83316668797592
$ws.Range("F24").Value = 38.67506756041122
$ws.Range("G24").Value = 44.98076677915725
$ws.Range("H24").Value = 18.02553497256426
$ws.Range("J24").Value = 10.04012856804
$ws.Range("K24").Value = 14.46342431338525
$ws.Range("M24").Value = 17.25633554061724
$ws.Range("N24").Value = 20.33469078380329
$ws.Range("B25").Value = 7.936988797606957
$ws.Range("D25").Value = 7.821018969271554
$ws.Range("E25").Value = 12.83860492807831
$ws.Range("F25").Value = 38.58619003344756
$ws.Range("G25").Value = 44.70952831655131
$ws.Range("H25").Value = 18.06820608415866
$ws.Range("J25").Value = 10.07432021246648
$ws.Range("K25").Value = 13.94472136185739
$ws.Range("M25").Value = 17.1226915378806
$ws.Range("N25").Value = 20.4551462736904
